$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new rows above row 519 (old rows 519.. shift down by 2,
# ending up at 521..564). Insert twice at the same index so the second
# insert pushes the first blank row further down too.
$ws.Rows.Item(519).Insert()
$ws.Rows.Item(519).Insert()

# New row 519 data
$ws.Cells.Item(519,1).Value  = 6
$ws.Cells.Item(519,2).Value  = "Mercado Mayorista Lo Valledor de Santiago"
$ws.Cells.Item(519,3).Value  = "Metropolitana"
$ws.Cells.Item(519,4).Value  = 44783
$ws.Cells.Item(519,5).Value  = 13
$ws.Cells.Item(519,6).Value  = 100112044
$ws.Cells.Item(519,7).Value  = "Perejil"
$ws.Cells.Item(519,8).Value  = "Sin especificar"
$ws.Cells.Item(519,9).Value  = "Primera"
$ws.Cells.Item(519,10).Value = 110
$ws.Cells.Item(519,11).Value = 20000
$ws.Cells.Item(519,12).Value = 20000
$ws.Cells.Item(519,13).Value = 20000
$ws.Cells.Item(519,14).Value = "`$/docena de atados"
$ws.Cells.Item(519,15).Value = "Región Metropolitana"
$ws.Cells.Item(519,16).Value = 6667
$ws.Cells.Item(519,17).Value = 3
$ws.Cells.Item(519,18).Value = "Hortaliza"

# New row 520 data
$ws.Cells.Item(520,1).Value  = 6
$ws.Cells.Item(520,2).Value  = "Mercado Mayorista Lo Valledor de Santiago"
$ws.Cells.Item(520,3).Value  = "Metropolitana"
$ws.Cells.Item(520,4).Value  = 44783
$ws.Cells.Item(520,5).Value  = 13
$ws.Cells.Item(520,6).Value  = 100112044
$ws.Cells.Item(520,7).Value  = "Perejil"
$ws.Cells.Item(520,8).Value  = "Sin especificar"
$ws.Cells.Item(520,9).Value  = "Segunda"
$ws.Cells.Item(520,10).Value = 60
$ws.Cells.Item(520,11).Value = 17000
$ws.Cells.Item(520,12).Value = 17000
$ws.Cells.Item(520,13).Value = 17000
$ws.Cells.Item(520,14).Value = "`$/docena de atados"
$ws.Cells.Item(520,15).Value = "Región Metropolitana"
$ws.Cells.Item(520,16).Value = 5667
$ws.Cells.Item(520,17).Value = 3
$ws.Cells.Item(520,18).Value = "Hortaliza"

# Make sure the date cells carry the same number format style as the
# other "Fecha" column cells (column D uses style index referencing a
# date/time number format).
$ws.Cells.Item(519,4).NumberFormat = $ws.Cells.Item(521,4).NumberFormat
$ws.Cells.Item(520,4).NumberFormat = $ws.Cells.Item(521,4).NumberFormat
